$wb = $excel.ActiveWorkbook

# --- BFS ---
$ws = $wb.Worksheets.Item("BFS")
$ws.Range("B2").Value = 10
$ws.Range("D2").Value = "[10, 14, 18, 19]"
$ws.Range("E2").Value = 390
$ws.Range("F2").Value = 15
$ws.Range("G2").Value = 1.333333333333333
$ws.Range("H2").Value = 0.0002052783966064453
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = "[25, 3, 6, 9, 13, 17, 20, 24]"
$ws.Range("E3").Value = 1153
$ws.Range("F3").Value = 31
$ws.Range("G3").Value = 1.096774193548387
$ws.Range("H3").Value = 0.0000922679901123046875
$ws.Range("B4").Value = 31
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = "[31, 38, 39, 40]"
$ws.Range("E4").Value = 204
$ws.Range("F4").Value = 21
$ws.Range("G4").Value = 1.333333333333333
$ws.Range("H4").Value = 0.0001003742218017578
$ws.Range("B5").Value = 41
$ws.Range("C5").Value = 29
$ws.Range("D5").Value = "[41, 38, 31, 30, 29]"
$ws.Range("E5").Value = 355
$ws.Range("F5").Value = 16
$ws.Range("G5").Value = 1.25
$ws.Range("H5").Value = 0.0001046657562255859
$ws.Range("B6").Value = 14
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = "[14, 15, 16, 17, 30]"
$ws.Range("E6").Value = 401
$ws.Range("F6").Value = 20
$ws.Range("G6").Value = 1.3
$ws.Range("H6").Value = 0.0001049041748046875

# --- DFS ---
$ws = $wb.Worksheets.Item("DFS")
$ws.Range("B2").Value = 10
$ws.Range("D2").Value = "[10, 14, 15, 16, 19]"
$ws.Range("E2").Value = 380
$ws.Range("F2").Value = 3121139
$ws.Range("G2").Value = 0.9999996796041445
$ws.Range("H2").Value = 1.29628849029541
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = "[25, 26, 27, 28, 29, 30, 31, 32, 24]"
$ws.Range("E3").Value = 615
$ws.Range("F3").Value = 2716923
$ws.Range("G3").Value = 0.9999996319365694
$ws.Range("H3").Value = 1.133981943130493
$ws.Range("B4").Value = 31
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = "[31, 38, 39, 40]"
$ws.Range("E4").Value = 204
$ws.Range("F4").Value = 6634523
$ws.Range("G4").Value = 0.9999998492732635
$ws.Range("H4").Value = 2.81198787689209
$ws.Range("B5").Value = 41
$ws.Range("C5").Value = 29
$ws.Range("D5").Value = "[41, 38, 37, 30, 29]"
$ws.Range("E5").Value = 352
$ws.Range("F5").Value = 3349902
$ws.Range("G5").Value = 0.9999997014838046
$ws.Range("H5").Value = 1.245441436767578
$ws.Range("B6").Value = 14
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = "[14, 15, 16, 17, 30]"
$ws.Range("E6").Value = 401
$ws.Range("F6").Value = 2153983
$ws.Range("G6").Value = 0.9999995357437825
$ws.Range("H6").Value = 0.878685712814331

# --- BCU ---
$ws = $wb.Worksheets.Item("BCU")
$ws.Range("B2").Value = 10
$ws.Range("D2").Value = "[10, 14, 15, 16, 19]"
$ws.Range("E2").Value = 380
$ws.Range("F2").Value = 14
$ws.Range("G2").Value = 0.6190476190476191
$ws.Range("H2").Value = 0.0001392364501953125
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = "[25, 26, 27, 28, 29, 30, 31, 32, 24]"
$ws.Range("E3").Value = 615
$ws.Range("F3").Value = 35
$ws.Range("G3").Value = 0.7906976744186046
$ws.Range("H3").Value = 0.0001776218414306641
$ws.Range("B4").Value = 31
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = "[31, 38, 39, 40]"
$ws.Range("E4").Value = 204
$ws.Range("F4").Value = 14
$ws.Range("G4").Value = 0.5909090909090909
$ws.Range("H4").Value = 0.0000894069671630859375
$ws.Range("B5").Value = 41
$ws.Range("C5").Value = 29
$ws.Range("D5").Value = "[41, 38, 37, 30, 29]"
$ws.Range("E5").Value = 352
$ws.Range("F5").Value = 18
$ws.Range("G5").Value = 0.7727272727272727
$ws.Range("H5").Value = 0.0000884532928466796875
$ws.Range("B6").Value = 14
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = "[14, 15, 16, 17, 30]"
$ws.Range("E6").Value = 401
$ws.Range("F6").Value = 20
$ws.Range("G6").Value = 0.76
$ws.Range("H6").Value = 0.0001554489135742188

# --- A_Estrela_Euclidiano ---
$ws = $wb.Worksheets.Item("A_Estrela_Euclidiano")
$ws.Range("B2").Value = 10
$ws.Range("D2").Value = "[10, 14, 15, 16, 19]"
$ws.Range("E2").Value = 380
$ws.Range("F2").Value = 13
$ws.Range("G2").Value = 3.071428571428572
$ws.Range("H2").Value = 0.0001668930053710938
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = "[25, 26, 27, 28, 29, 30, 31, 32, 24]"
$ws.Range("E3").Value = 615
$ws.Range("F3").Value = 32
$ws.Range("G3").Value = 3.071428571428572
$ws.Range("H3").Value = 0.0001673698425292969
$ws.Range("B4").Value = 31
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = "[31, 38, 39, 40]"
$ws.Range("E4").Value = 204
$ws.Range("F4").Value = 12
$ws.Range("G4").Value = 3.071428571428572
$ws.Range("H4").Value = 0.0001206398010253906
$ws.Range("B5").Value = 41
$ws.Range("C5").Value = 29
$ws.Range("D5").Value = "[41, 38, 37, 30, 29]"
$ws.Range("E5").Value = 352
$ws.Range("F5").Value = 15
$ws.Range("G5").Value = 3.071428571428572
$ws.Range("H5").Value = 0.0001127719879150391
$ws.Range("B6").Value = 14
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = "[14, 15, 16, 17, 30]"
$ws.Range("E6").Value = 401
$ws.Range("G6").Value = 3.071428571428572
$ws.Range("H6").Value = 0.0001287460327148438

# --- A_Estrela_Haversiano ---
$ws = $wb.Worksheets.Item("A_Estrela_Haversiano")
$ws.Range("B2").Value = 10
$ws.Range("D2").Value = "[10, 14, 15, 16, 19]"
$ws.Range("E2").Value = 380
$ws.Range("F2").Value = 13
$ws.Range("G2").Value = 3.071428571428572
$ws.Range("H2").Value = 0.0002303123474121094
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = "[25, 26, 27, 28, 29, 30, 31, 32, 24]"
$ws.Range("E3").Value = 615
$ws.Range("F3").Value = 32
$ws.Range("G3").Value = 3.071428571428572
$ws.Range("H3").Value = 0.0001883506774902344
$ws.Range("B4").Value = 31
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = "[31, 38, 39, 40]"
$ws.Range("E4").Value = 204
$ws.Range("F4").Value = 12
$ws.Range("G4").Value = 3.071428571428572
$ws.Range("H4").Value = 0.0001289844512939453
$ws.Range("B5").Value = 41
$ws.Range("C5").Value = 29
$ws.Range("D5").Value = "[41, 38, 37, 30, 29]"
$ws.Range("E5").Value = 352
$ws.Range("F5").Value = 15
$ws.Range("G5").Value = 3.071428571428572
$ws.Range("H5").Value = 0.0001378059387207031
$ws.Range("B6").Value = 14
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = "[14, 15, 16, 17, 30]"
$ws.Range("E6").Value = 401
$ws.Range("G6").Value = 3.071428571428572
$ws.Range("H6").Value = 0.0001387596130371094
